$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 720.5
$ws.Cells.Item(8, 9).Value = 720.5
$ws.Cells.Item(8, 11).Value = 2161.5
$ws.Cells.Item(8, 13).Value = -2022.5
$ws.Cells.Item(64, 8).Value = 8782.593999999999
$ws.Cells.Item(64, 9).Value = 5002.846
$ws.Cells.Item(64, 10).Value = 11368.737
$ws.Cells.Item(64, 11).Value = 5002.846
$ws.Cells.Item(64, 12).Value = 11368.737
$ws.Cells.Item(64, 13).Value = -4754.846
$ws.Cells.Item(64, 14).Value = -11864.737
$ws.Cells.Item(67, 8).Value = 8782.593999999999
$ws.Cells.Item(67, 9).Value = 5002.846
$ws.Cells.Item(67, 10).Value = 11368.737
$ws.Cells.Item(67, 11).Value = 5002.846
$ws.Cells.Item(67, 12).Value = 11368.737
$ws.Cells.Item(67, 13).Value = -4144.846
$ws.Cells.Item(67, 14).Value = -13084.737
$ws.Cells.Item(113, 8).Value = 4700
$ws.Cells.Item(113, 9).Value = 3900
$ws.Cells.Item(113, 10).Value = 4966.6665
$ws.Cells.Item(113, 11).Value = 3900
$ws.Cells.Item(113, 12).Value = 4966.6665
$ws.Cells.Item(113, 13).Value = -646
$ws.Cells.Item(113, 14).Value = -11474.6665
$ws.Cells.Item(129, 8).Value = 1740.0625
$ws.Cells.Item(129, 10).Value = 2337.5
$ws.Cells.Item(129, 12).Value = 7012.5
$ws.Cells.Item(129, 14).Value = -17012.5
$ws.Cells.Item(135, 8).Value = 1795.4286
$ws.Cells.Item(135, 9).Value = 1795.4286
$ws.Cells.Item(135, 11).Value = 16158.8574
$ws.Cells.Item(135, 13).Value = -13623.8574
$ws.Cells.Item(137, 8).Value = 2254.8928
$ws.Cells.Item(137, 9).Value = 1992.3334
$ws.Cells.Item(137, 10).Value = 3042.5715
$ws.Cells.Item(137, 11).Value = 5977.0002
$ws.Cells.Item(137, 12).Value = 9127.7145
$ws.Cells.Item(137, 13).Value = -3427.0002
$ws.Cells.Item(137, 14).Value = -14227.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5023.9
$ws.Cells.Item(2, 9).Value = 4656
$ws.Cells.Item(2, 10).Value = 5882.3335
$ws.Cells.Item(2, 11).Value = 4656
$ws.Cells.Item(2, 12).Value = 5882.3335
$ws.Cells.Item(2, 13).Value = -4543
$ws.Cells.Item(2, 14).Value = -6108.3335
$ws.Cells.Item(24, 8).Value = 14000
$ws.Cells.Item(24, 10).Value = 14000
$ws.Cells.Item(24, 12).Value = 14000
$ws.Cells.Item(24, 14).Value = -14748
$ws.Cells.Item(32, 8).Value = 2594.353
$ws.Cells.Item(32, 9).Value = 2485.6875
$ws.Cells.Item(32, 11).Value = 2485.6875
$ws.Cells.Item(32, 13).Value = -2198.6875
$ws.Cells.Item(61, 8).Value = 5699.75
$ws.Cells.Item(61, 9).Value = 5926.467
$ws.Cells.Item(61, 10).Value = 2299
$ws.Cells.Item(61, 11).Value = 5926.467
$ws.Cells.Item(61, 12).Value = 2299
$ws.Cells.Item(61, 13).Value = -5714.467
$ws.Cells.Item(61, 14).Value = -2723
$ws.Cells.Item(74, 8).Value = 3673.3845
$ws.Cells.Item(74, 9).Value = 3078.5
$ws.Cells.Item(74, 10).Value = 5656.3335
$ws.Cells.Item(74, 11).Value = 3078.5
$ws.Cells.Item(74, 12).Value = 5656.3335
$ws.Cells.Item(74, 13).Value = -2204.5
$ws.Cells.Item(74, 14).Value = -7404.3335
$ws.Cells.Item(77, 8).Value = 3673.3845
$ws.Cells.Item(77, 9).Value = 3078.5
$ws.Cells.Item(77, 10).Value = 5656.3335
$ws.Cells.Item(77, 11).Value = 15392.5
$ws.Cells.Item(77, 12).Value = 28281.6675
$ws.Cells.Item(77, 13).Value = -11024.5
$ws.Cells.Item(77, 14).Value = -37017.6675
$ws.Cells.Item(94, 8).Value = 25400
$ws.Cells.Item(94, 10).Value = 25400
$ws.Cells.Item(94, 12).Value = 25400
$ws.Cells.Item(94, 14).Value = -27202
$ws.Cells.Item(100, 8).Value = 14000
$ws.Cells.Item(100, 10).Value = 14000
$ws.Cells.Item(100, 12).Value = 14000
$ws.Cells.Item(100, 14).Value = -16164
$ws.Cells.Item(101, 8).Value = 40000
$ws.Cells.Item(101, 10).Value = 40000
$ws.Cells.Item(101, 12).Value = 40000
$ws.Cells.Item(101, 14).Value = -46490
$ws.Cells.Item(102, 8).Value = 2097.2144
$ws.Cells.Item(102, 9).Value = 1804.5186
$ws.Cells.Item(102, 11).Value = 1804.5186
$ws.Cells.Item(102, 13).Value = -182.5186000000001
$ws.Cells.Item(116, 8).Value = 5023.9
$ws.Cells.Item(116, 9).Value = 4656
$ws.Cells.Item(116, 10).Value = 5882.3335
$ws.Cells.Item(116, 11).Value = 4656
$ws.Cells.Item(116, 12).Value = 5882.3335
$ws.Cells.Item(116, 13).Value = -2362
$ws.Cells.Item(116, 14).Value = -10470.3335
$ws.Cells.Item(132, 8).Value = 2695.875
$ws.Cells.Item(132, 9).Value = 2408.9333
$ws.Cells.Item(132, 11).Value = 7226.7999
$ws.Cells.Item(132, 13).Value = -4696.7999
$ws.Cells.Item(136, 8).Value = 5699.75
$ws.Cells.Item(136, 9).Value = 5926.467
$ws.Cells.Item(136, 10).Value = 2299
$ws.Cells.Item(136, 11).Value = 17779.401
$ws.Cells.Item(136, 12).Value = 6897
$ws.Cells.Item(136, 13).Value = -15229.401
$ws.Cells.Item(136, 14).Value = -11997

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5023.9
$ws.Cells.Item(3, 9).Value = 4656
$ws.Cells.Item(3, 10).Value = 5882.3335
$ws.Cells.Item(3, 11).Value = 4656
$ws.Cells.Item(3, 12).Value = 5882.3335
$ws.Cells.Item(3, 13).Value = -4542
$ws.Cells.Item(3, 14).Value = -6110.3335
$ws.Cells.Item(20, 8).Value = 1557.6842
$ws.Cells.Item(20, 9).Value = 1499.7858
$ws.Cells.Item(20, 11).Value = 1499.7858
$ws.Cells.Item(20, 13).Value = -1252.7858
$ws.Cells.Item(105, 8).Value = 2005.6875
$ws.Cells.Item(105, 9).Value = 1989.742
$ws.Cells.Item(105, 11).Value = 1989.742
$ws.Cells.Item(105, 13).Value = -242.742
$ws.Cells.Item(132, 8).Value = 64998.668
$ws.Cells.Item(132, 10).Value = 64998.668
$ws.Cells.Item(132, 12).Value = 64998.668
$ws.Cells.Item(132, 14).Value = -75118.66800000001
$ws.Cells.Item(134, 8).Value = 2832.5
$ws.Cells.Item(134, 9).Value = 2856.6667
$ws.Cells.Item(134, 11).Value = 8570.000100000001
$ws.Cells.Item(134, 13).Value = -6035.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3630.0588
$ws.Cells.Item(31, 9).Value = 1794.4814
$ws.Cells.Item(31, 10).Value = 10710.143
$ws.Cells.Item(31, 11).Value = 1794.4814
$ws.Cells.Item(31, 12).Value = 10710.143
$ws.Cells.Item(31, 13).Value = -1499.4814
$ws.Cells.Item(31, 14).Value = -11300.143
$ws.Cells.Item(34, 8).Value = 3630.0588
$ws.Cells.Item(34, 9).Value = 1794.4814
$ws.Cells.Item(34, 10).Value = 10710.143
$ws.Cells.Item(34, 11).Value = 1794.4814
$ws.Cells.Item(34, 12).Value = 10710.143
$ws.Cells.Item(34, 13).Value = -1592.4814
$ws.Cells.Item(34, 14).Value = -11114.143
$ws.Cells.Item(95, 8).Value = 12779.8
$ws.Cells.Item(95, 10).Value = 12779.8
$ws.Cells.Item(95, 12).Value = 12779.8
$ws.Cells.Item(95, 14).Value = -18271.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 2344.25
$ws.Cells.Item(107, 10).Value = 933
$ws.Cells.Item(107, 12).Value = 2799
$ws.Cells.Item(107, 14).Value = -6639
$ws.Cells.Item(122, 8).Value = 3025.6667
$ws.Cells.Item(122, 10).Value = 3995
$ws.Cells.Item(122, 12).Value = 35955
$ws.Cells.Item(122, 14).Value = -40855
$ws.Cells.Item(132, 8).Value = 2639.5454
$ws.Cells.Item(132, 9).Value = 1678.6666
$ws.Cells.Item(132, 11).Value = 15107.9994
$ws.Cells.Item(132, 13).Value = -12577.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(22, 8).Value = 9134.666999999999
$ws.Cells.Item(22, 9).Value = 6800
$ws.Cells.Item(22, 10).Value = 11469.333
$ws.Cells.Item(22, 11).Value = 6800
$ws.Cells.Item(22, 12).Value = 11469.333
$ws.Cells.Item(22, 13).Value = -6271
$ws.Cells.Item(22, 14).Value = -12527.333
$ws.Cells.Item(70, 8).Value = 23815882
$ws.Cells.Item(70, 9).Value = 83337336
$ws.Cells.Item(70, 10).Value = 7300.8
$ws.Cells.Item(70, 11).Value = 83337336
$ws.Cells.Item(70, 12).Value = 7300.8
$ws.Cells.Item(70, 13).Value = -83337066
$ws.Cells.Item(70, 14).Value = -7840.8
$ws.Cells.Item(73, 8).Value = 23815882
$ws.Cells.Item(73, 9).Value = 83337336
$ws.Cells.Item(73, 10).Value = 7300.8
$ws.Cells.Item(73, 11).Value = 83337336
$ws.Cells.Item(73, 12).Value = 7300.8
$ws.Cells.Item(73, 13).Value = -83336400
$ws.Cells.Item(73, 14).Value = -9172.799999999999
$ws.Cells.Item(132, 8).Value = 2284.6956
$ws.Cells.Item(132, 9).Value = 2384.6191
$ws.Cells.Item(132, 10).Value = 1235.5
$ws.Cells.Item(132, 11).Value = 7153.8573
$ws.Cells.Item(132, 12).Value = 3706.5
$ws.Cells.Item(132, 13).Value = -4623.8573
$ws.Cells.Item(132, 14).Value = -8766.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 795.5
$ws.Cells.Item(16, 9).Value = 554.8
$ws.Cells.Item(16, 11).Value = 554.8
$ws.Cells.Item(16, 13).Value = -384.8
$ws.Cells.Item(132, 8).Value = 5677.1304
$ws.Cells.Item(132, 9).Value = 5844
$ws.Cells.Item(132, 11).Value = 17532
$ws.Cells.Item(132, 13).Value = -15002

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(101, 8).Value = 19000
$ws.Cells.Item(101, 10).Value = 19000
$ws.Cells.Item(101, 12).Value = 19000
$ws.Cells.Item(101, 14).Value = -25490
$ws.Cells.Item(136, 8).Value = 1465.8235
$ws.Cells.Item(136, 9).Value = 1344.6786
$ws.Cells.Item(136, 10).Value = 2031.1666
$ws.Cells.Item(136, 11).Value = 4034.0358
$ws.Cells.Item(136, 12).Value = 6093.4998
$ws.Cells.Item(136, 13).Value = -1484.0358
$ws.Cells.Item(136, 14).Value = -11193.4998
